$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.455.14'
Set-TextValue $ws.Range("E2") '  -1.26%  '
Set-TextValue $ws.Range("D3") '1.849.45'
Set-TextValue $ws.Range("E3") '  -0.73%  '
Set-TextValue $ws.Range("D4") '0.9993'
Set-TextValue $ws.Range("E4") '  -0.18%  '
Set-TextValue $ws.Range("D5") '265.36'
Set-TextValue $ws.Range("E5") '  -1.81%  '
Set-TextValue $ws.Range("D6") '0.9996'
Set-TextValue $ws.Range("E6") '  -0.05%  '
Set-TextValue $ws.Range("D7") '0.5216'
Set-TextValue $ws.Range("E7") '  -1.61%  '
Set-TextValue $ws.Range("D8") '0.3277'
Set-TextValue $ws.Range("E8") '  -2.42%  '
Set-TextValue $ws.Range("D9") '0.06825'
Set-TextValue $ws.Range("E9") '  +0.41%  '
Set-TextValue $ws.Range("D10") '18.90'
Set-TextValue $ws.Range("E10") '  -4.08%  '
Set-TextValue $ws.Range("D11") '0.7787'
Set-TextValue $ws.Range("E11") '  -1.05%  '
Set-TextValue $ws.Range("D12") '0.07784'
Set-TextValue $ws.Range("E12") '  +0.38%  '
Set-TextValue $ws.Range("D13") '1.813.69'
Set-TextValue $ws.Range("E13") '  -2.64%  '
Set-TextValue $ws.Range("D14") '88.30'
Set-TextValue $ws.Range("E14") '  -1.69%  '
Set-TextValue $ws.Range("D15") '5.026'
Set-TextValue $ws.Range("E15") '  -1.60%  '
Set-TextValue $ws.Range("E16") '  +0.01%  '
Set-TextValue $ws.Range("D17") '13.97'
Set-TextValue $ws.Range("E17") '  -2.91%  '
Set-TextValue $ws.Range("D18") '0.000008001'
Set-TextValue $ws.Range("E18") '  -0.01%  '
Set-TextValue $ws.Range("E19") '  +0.06%  '
Set-TextValue $ws.Range("D20") '26.438.80'
Set-TextValue $ws.Range("E20") '  -1.38%  '
Set-TextValue $ws.Range("D21") '2.066.38'
Set-TextValue $ws.Range("E21") '  -2.09%  '
Set-TextValue $ws.Range("D22") '4.646'
Set-TextValue $ws.Range("E22") '  +0.02%  '
Set-TextValue $ws.Range("D23") '9.585'
Set-TextValue $ws.Range("E23") '  -3.12%  '
Set-TextValue $ws.Range("D24") '6.012'
Set-TextValue $ws.Range("E24") '  -0.71%  '
Set-TextValue $ws.Range("D25") '144.46'
Set-TextValue $ws.Range("E25") '  -0.96%  '
Set-TextValue $ws.Range("D26") '2.196'
Set-TextValue $ws.Range("E26") '  -8.09%  '
Set-TextValue $ws.Range("D27") '1.669'
Set-TextValue $ws.Range("E27") '  +1.05%  '
Set-TextValue $ws.Range("D28") '17.05'
Set-TextValue $ws.Range("E28") '  -0.73%  '
Set-TextValue $ws.Range("D29") '112.32'
Set-TextValue $ws.Range("E29") '  -0.47%  '
Set-TextValue $ws.Range("D30") '4.185'
Set-TextValue $ws.Range("E30") '  -2.98%  '
Set-TextValue $ws.Range("D31") '4.151'
Set-TextValue $ws.Range("E31") '  -2.81%  '
Set-TextValue $ws.Range("D32") '0.08768'
Set-TextValue $ws.Range("E32") '  -0.97%  '
Set-TextValue $ws.Range("D33") '0.04845'
Set-TextValue $ws.Range("E33") '  -2.08%  '
Set-TextValue $ws.Range("D34") '1.139'
Set-TextValue $ws.Range("E34") '  -1.35%  '
Set-TextValue $ws.Range("D35") '0.7213'
Set-TextValue $ws.Range("E35") '  -0.33%  '
Set-TextValue $ws.Range("D36") '2.860'
Set-TextValue $ws.Range("E36") '  -0.72%  '
Set-TextValue $ws.Range("D37") '3.097'
Set-TextValue $ws.Range("E37") '  -3.04%  '
Set-TextValue $ws.Range("D38") '0.01782'
Set-TextValue $ws.Range("E38") '  -2.87%  '
Set-TextValue $ws.Range("D39") '2.220'
Set-TextValue $ws.Range("E39") '  -3.43%  '
Set-TextValue $ws.Range("D40") '0.4886'
Set-TextValue $ws.Range("E40") '  -3.45%  '
Set-TextValue $ws.Range("D41") '0.9092'
Set-TextValue $ws.Range("E41") '  -1.27%  '
Set-TextValue $ws.Range("D42") '111.48'
Set-TextValue $ws.Range("E42") '  -3.79%  '
Set-TextValue $ws.Range("D43") '6.076'
Set-TextValue $ws.Range("E43") '  -0.75%  '
Set-TextValue $ws.Range("D44") '1.001'
Set-TextValue $ws.Range("E44") '  +0.10%  '
Set-TextValue $ws.Range("D45") '7.736'
Set-TextValue $ws.Range("E45") '  -2.71%  '
Set-TextValue $ws.Range("D48") '9.135'
Set-TextValue $ws.Range("E48") '  -1.98%  '
Set-TextValue $ws.Range("D49") '0.1239'
Set-TextValue $ws.Range("E49") '  -6.19%  '
Set-TextValue $ws.Range("D50") '35.03'
Set-TextValue $ws.Range("E50") '  -2.43%  '
Set-TextValue $ws.Range("D51") '0.8942'
Set-TextValue $ws.Range("E51") '  +2.32%  '

# Rows 46 and 47: coin identities swap with new values
Set-TextValue $ws.Range("B46") 'Decentraland'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range("D46") '0.4178'
Set-TextValue $ws.Range("E46") '  -4.65%  '

Set-TextValue $ws.Range("B47") 'Cronos'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D47") '0.05944'
Set-TextValue $ws.Range("E47") '  +0.13%  '
